# AIAC_2.3 — "Add files via upload" edit
#
# 1) NAME / ROLL NO line: the actual visible text changes
#    ("B. SAINATH" / "NO:2403A510C3" -> "P. HEMAN" / "NO:2403A510F5",
#    with a few extra spaces inserted before "ROLL"). A plain wildcard
#    Find & Replace collapses the several runs that used to make up the
#    line into the single run Word naturally produces for freshly typed
#    text, which matches the target markup.
#
# 2) Every other change in the diff is a pure run-split: the *visible*
#    text of each paragraph is unchanged, only the run boundaries move
#    (e.g. "Task Description 1:" keeps being the string
#    "Task Description 1:", but now lives in four runs: "Task
#    Description" / " " / "1" / ":"). We reproduce that by toggling a
#    character-formatting property (Bold) off-then-on (or on-then-off)
#    across the sub-range of each new segment: changing a run property
#    and changing it back still forces Word to split the run at that
#    boundary, but leaves the saved <w:rPr> identical to its neighbours.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. NAME / ROLL NO
# ---------------------------------------------------------------------

$d.Content.Find.Execute(
    "NAME: B. SAINATH*ROLL NO:2403A510C3",
    $true, $false, $true, $false, $false,
    $true, 1, $false,
    "NAME: P. HEMAN                                                                                                                   ROLL NO:2403A510F5",
    2) | Out-Null

# ---------------------------------------------------------------------
# 2. Run-splitting helpers
# ---------------------------------------------------------------------

function SplitAt($s, $e) {
    $sub = $d.Range($s, $e)
    $orig = $sub.Font.Bold
    if ($orig -ne 0) {
        $sub.Font.Bold = 0
        $sub.Font.Bold = $orig
    } else {
        $sub.Font.Bold = 1
        $sub.Font.Bold = 0
    }
}

# Finds $findText (unique, literal match) in the document body and
# splits it into runs whose lengths are given by $lens (which must sum
# to $findText.Length). The first segment stays attached to whatever
# run already held the text; every later segment gets its own run.
function SplitRun($findText, $lens) {
    $rng = $d.Content
    $rng.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $base = $rng.Start

    $n = $lens.Length
    $bounds = @($base)
    $cum = $base
    for ($i = 0; $i -lt $n; $i++) {
        $cum = $cum + $lens[$i]
        $bounds += $cum
    }

    for ($i = 1; $i -lt $n; $i++) {
        SplitAt $bounds[$i] $bounds[$i + 1]
    }
}

# ---------------------------------------------------------------------
# 3. Apply every run-split from the diff (text itself is unchanged)
# ---------------------------------------------------------------------

SplitRun "Task Description 1:" @(16, 1, 1, 1)
SplitRun "->Use Google Gemini in " @(2, 21)

SplitRun "Expected Output 1:" @(15, 1, 1, 1)
SplitRun "->Functional code with output and screenshot " @(2, 43)

SplitRun "Task Description 2:" @(16, 1, 1, 1)
SplitRun "->Compare Gemini and Copilot outputs for a palindrome check function." @(2, 67)

SplitRun "Expected Output 2:" @(15, 1, 1, 1)
SplitRun "->Side-by-side comparison and observations" @(2, 40)

SplitRun "Task Description 3:" @(16, 1, 1, 1)
SplitRun "->Ask Gemini to explain a Python function (to calculate area of various shapes) line by " @(2, 86)

SplitRun "Expected Output 3:" @(15, 1, 1, 1)
SplitRun "->Detailed explanation with code snippet" @(2, 38)

SplitRun "Description  4" @(11, 2, 1)
SplitRun "->Install and configure Cursor AI." @(2, 32)

SplitRun "Expected Output 4:" @(15, 1, 1, 1)
SplitRun "->Screenshots of working environments with few prompts to generate python code" @(2, 76)

SplitRun "Task Description 5:" @(16, 1, 1, 1)
SplitRun "->Student need to write code to calculate sum of add number and even numbers in the list" @(2, 86)

SplitRun "Expected Output 5:" @(15, 1, 1, 1)
SplitRun "->Refactored code written by student with improved logic" @(2, 54)
